# Weekly update: add a new week (Región del Maule) on top of the
# "Terminal Hortofrutícola Agro Chillán - Pimiento" Zafiro rojo/verde
# block, pushing the previously-newest rows (253-256) down to 255-258.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two fresh rows above the current row 253; everything that was
# at 253+ (including its row formatting) shifts down to 255+.
$ws.Rows("253:254").Insert()

# --- Row 253: Zafiro rojo, new week (2022-04-05, serial 44656) ---
$ws.Cells.Item(253, 1).Value = 7
$ws.Cells.Item(253, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(253, 3).Value = "Ñuble"
$ws.Cells.Item(253, 4).Value = 44656
$ws.Cells.Item(253, 5).Value = 16
$ws.Cells.Item(253, 6).Value = 100112002
$ws.Cells.Item(253, 7).Value = "Pimiento"
$ws.Cells.Item(253, 8).Value = "Zafiro rojo"
$ws.Cells.Item(253, 9).Value = "Primera"
$ws.Cells.Item(253, 10).Value = 120
$ws.Cells.Item(253, 11).Value = 15000
$ws.Cells.Item(253, 12).Value = 16000
$ws.Cells.Item(253, 13).Value = 15500
$ws.Cells.Item(253, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(253, 15).Value = "Región del Maule"
$ws.Cells.Item(253, 16).Value = 1033
$ws.Cells.Item(253, 17).Value = 15
$ws.Cells.Item(253, 18).Value = "Hortaliza"

# --- Row 254: Zafiro verde, new week (2022-04-05, serial 44656) ---
$ws.Cells.Item(254, 1).Value = 7
$ws.Cells.Item(254, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(254, 3).Value = "Ñuble"
$ws.Cells.Item(254, 4).Value = 44656
$ws.Cells.Item(254, 5).Value = 16
$ws.Cells.Item(254, 6).Value = 100112002
$ws.Cells.Item(254, 7).Value = "Pimiento"
$ws.Cells.Item(254, 8).Value = "Zafiro verde"
$ws.Cells.Item(254, 9).Value = "Primera"
$ws.Cells.Item(254, 10).Value = 120
$ws.Cells.Item(254, 11).Value = 9000
$ws.Cells.Item(254, 12).Value = 10000
$ws.Cells.Item(254, 13).Value = 9500
$ws.Cells.Item(254, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(254, 15).Value = "Región del Maule"
$ws.Cells.Item(254, 16).Value = 633
$ws.Cells.Item(254, 17).Value = 15
$ws.Cells.Item(254, 18).Value = "Hortaliza"
